# Insert a new data row before the current row 82 (pushes existing rows 82-185
# down to 83-186), then populate the newly inserted row 82 with the new
# weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 82, shifting rows 82:185 down to 83:186
$ws.Rows(82).Insert()

# Populate the new row 82 with the new record's data
$ws.Range("A82").Value = 6
$ws.Range("B82").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C82").Value = "Metropolitana"
$ws.Range("D82").Value = 44638
$ws.Range("E82").Value = 13
$ws.Range("F82").Value = "Fruta"
$ws.Range("G82").Value = 100101
$ws.Range("H82").Value = "Berries"
$ws.Range("I82").Value = 100101004
$ws.Range("J82").Value = "Frambuesa"
$ws.Range("K82").Value = "Sin especificar"
$ws.Range("L82").Value = "Especial"
$ws.Range("M82").Value = 300
$ws.Range("N82").Value = 8000
$ws.Range("O82").Value = 8000
$ws.Range("P82").Value = 8000
$ws.Range("Q82").Value = "$/bandeja 2 kilos"
$ws.Range("R82").Value = "Provincia de Linares"
$ws.Range("S82").Value = 4000
$ws.Range("T82").Value = 2
